# Calibration data cleanup: the raw curvature samples (rows 2-18, columns
# A:D) were captured out of chronological order ("time (s)" in column A).
# Sort the data block ascending by column A, leaving the header row intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:D18")
$keyRange  = $ws.Range("A2:A18")

$dataRange.Sort($keyRange, 1)
